# Weekly price update: insert a new record for the current week at row 204,
# pushing the existing historical rows (204-310) down by one (205-311).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 204 (shifts rows 204:310 down to 205:311)
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row with this week's data point
$ws.Cells.Item(204, 1).Value  = 10
$ws.Cells.Item(204, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(204, 3).Value  = "La Araucanía"
$ws.Cells.Item(204, 4).Value2 = 45097
$ws.Cells.Item(204, 5).Value  = 9
$ws.Cells.Item(204, 6).Value  = 100114007
$ws.Cells.Item(204, 7).Value  = "Jengibre"
$ws.Cells.Item(204, 8).Value  = "Sin especificar"
$ws.Cells.Item(204, 9).Value  = "Primera"
$ws.Cells.Item(204, 10).Value = 12
$ws.Cells.Item(204, 11).Value = 24000
$ws.Cells.Item(204, 12).Value = 24000
$ws.Cells.Item(204, 13).Value = 24000
$ws.Cells.Item(204, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(204, 15).Value = "Perú"
$ws.Cells.Item(204, 16).Value = 1846
$ws.Cells.Item(204, 17).Value = 13
$ws.Cells.Item(204, 18).Value = "Hortaliza"
